$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 348
$ws.Range("I2").Value = 1078
$ws.Range("J2").Value = 4281
$ws.Range("K2").Value = 26
$ws.Range("L2").Value = 1200
$ws.Range("M2").Value = 69
$ws.Range("N2").Value = 733
$ws.Range("O2").Value = 4
$ws.Range("P2").Value = 16
$ws.Range("Q2").Value = 13
$ws.Range("R2").Value = 59
$ws.Range("S2").Value = 445
$ws.Range("T2").Value = 737
$ws.Range("U2").Value = 47
$ws.Range("V2").Value = 6380
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 6522
$ws.Range("Y2").Value = 11
$ws.Range("Z2").Value = 116
$ws.Range("AA2").Value = 46
